$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that needs to be bumped by
# one day (46061 -> 46062) for every data row (rows 2 through 47).
for ($r = 2; $r -le 47; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value = 46062
    }
}
